# Auto-applied updates to Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Reflects refreshed market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 281.25
$ws.Cells.Item(9, 9).Value = 90
$ws.Cells.Item(9, 11).Value = 90
$ws.Cells.Item(9, 13).Value = 79

$ws.Cells.Item(18, 8).Value = 422
$ws.Cells.Item(18, 9).Value = 236.66667
$ws.Cells.Item(18, 11).Value = 236.66667
$ws.Cells.Item(18, 13).Value = 47.33332999999999

$ws.Cells.Item(33, 8).Value = 144.33333
$ws.Cells.Item(33, 9).Value = 144.33333
$ws.Cells.Item(33, 11).Value = 144.33333
$ws.Cells.Item(33, 13).Value = 84.66667000000001

$ws.Cells.Item(129, 8).Value = 954.1795
$ws.Cells.Item(129, 9).Value = 400.33334
$ws.Cells.Item(129, 10).Value = 1200.3334
$ws.Cells.Item(129, 11).Value = 1201.00002
$ws.Cells.Item(129, 12).Value = 3601.0002
$ws.Cells.Item(129, 13).Value = 3798.99998
$ws.Cells.Item(129, 14).Value = -13601.0002

$ws.Cells.Item(132, 8).Value = 3271.1482
$ws.Cells.Item(132, 9).Value = 3358.5
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 10075.5
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -7545.5
$ws.Cells.Item(132, 14).Value = -8060

$ws.Cells.Item(137, 8).Value = 1638.4
$ws.Cells.Item(137, 9).Value = 1455.238
$ws.Cells.Item(137, 10).Value = 2600
$ws.Cells.Item(137, 11).Value = 4365.714
$ws.Cells.Item(137, 12).Value = 7800
$ws.Cells.Item(137, 13).Value = -1815.714
$ws.Cells.Item(137, 14).Value = -12900

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6540.5
$ws.Cells.Item(32, 9).Value = 4690.0684
$ws.Cells.Item(32, 10).Value = 47250
$ws.Cells.Item(32, 11).Value = 4690.0684
$ws.Cells.Item(32, 12).Value = 47250
$ws.Cells.Item(32, 13).Value = -4403.0684
$ws.Cells.Item(32, 14).Value = -47824

$ws.Cells.Item(61, 8).Value = 1866.8387
$ws.Cells.Item(61, 9).Value = 1700.25
$ws.Cells.Item(61, 10).Value = 2438
$ws.Cells.Item(61, 11).Value = 1700.25
$ws.Cells.Item(61, 12).Value = 2438
$ws.Cells.Item(61, 13).Value = -1488.25
$ws.Cells.Item(61, 14).Value = -2862

$ws.Cells.Item(74, 8).Value = 62501020
$ws.Cells.Item(74, 9).Value = 90909800
$ws.Cells.Item(74, 10).Value = 1700
$ws.Cells.Item(74, 11).Value = 90909800
$ws.Cells.Item(74, 12).Value = 1700
$ws.Cells.Item(74, 13).Value = -90908926
$ws.Cells.Item(74, 14).Value = -3448

$ws.Cells.Item(77, 8).Value = 62501020
$ws.Cells.Item(77, 9).Value = 90909800
$ws.Cells.Item(77, 10).Value = 1700
$ws.Cells.Item(77, 11).Value = 454549000
$ws.Cells.Item(77, 12).Value = 8500
$ws.Cells.Item(77, 13).Value = -454544632
$ws.Cells.Item(77, 14).Value = -17236

$ws.Cells.Item(88, 8).Value = 127416.375
$ws.Cells.Item(88, 9).Value = 2074
$ws.Cells.Item(88, 10).Value = 202621.8
$ws.Cells.Item(88, 11).Value = 2074
$ws.Cells.Item(88, 12).Value = 202621.8
$ws.Cells.Item(88, 13).Value = -1668
$ws.Cells.Item(88, 14).Value = -203433.8

$ws.Cells.Item(91, 8).Value = 127416.375
$ws.Cells.Item(91, 9).Value = 2074
$ws.Cells.Item(91, 10).Value = 202621.8
$ws.Cells.Item(91, 11).Value = 2074
$ws.Cells.Item(91, 12).Value = 202621.8
$ws.Cells.Item(91, 13).Value = -670
$ws.Cells.Item(91, 14).Value = -205429.8

$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 3094.2778
$ws.Cells.Item(122, 9).Value = 2684.6155
$ws.Cells.Item(122, 10).Value = 4159.4
$ws.Cells.Item(122, 11).Value = 8053.8465
$ws.Cells.Item(122, 12).Value = 12478.2
$ws.Cells.Item(122, 13).Value = -5603.8465
$ws.Cells.Item(122, 14).Value = -17378.2

$ws.Cells.Item(132, 8).Value = 12084.234
$ws.Cells.Item(132, 9).Value = 1272.2195
$ws.Cells.Item(132, 10).Value = 85966.336
$ws.Cells.Item(132, 11).Value = 3816.6585
$ws.Cells.Item(132, 12).Value = 257899.008
$ws.Cells.Item(132, 13).Value = -1286.6585
$ws.Cells.Item(132, 14).Value = -262959.008

$ws.Cells.Item(136, 8).Value = 1866.8387
$ws.Cells.Item(136, 9).Value = 1700.25
$ws.Cells.Item(136, 10).Value = 2438
$ws.Cells.Item(136, 11).Value = 5100.75
$ws.Cells.Item(136, 12).Value = 7314
$ws.Cells.Item(136, 13).Value = -2550.75
$ws.Cells.Item(136, 14).Value = -12414

$ws.Cells.Item(139, 8).Value = 49626
$ws.Cells.Item(139, 10).Value = 49626
$ws.Cells.Item(139, 12).Value = 49626
$ws.Cells.Item(139, 14).Value = -59906

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 30000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()

$ws.Cells.Item(86, 8).Value = 1782.0322
$ws.Cells.Item(86, 10).Value = 2730
$ws.Cells.Item(86, 12).Value = 2730
$ws.Cells.Item(86, 14).Value = -4976

$ws.Cells.Item(89, 8).Value = 1782.0322
$ws.Cells.Item(89, 10).Value = 2730
$ws.Cells.Item(89, 12).Value = 13650
$ws.Cells.Item(89, 14).Value = -24882

$ws.Cells.Item(99, 8).Value = 1436.9546
$ws.Cells.Item(99, 9).Value = 1118.875
$ws.Cells.Item(99, 11).Value = 1118.875
$ws.Cells.Item(99, 13).Value = 379.125

$ws.Cells.Item(105, 8).Value = 3444
$ws.Cells.Item(105, 9).Value = 3489.1667
$ws.Cells.Item(105, 11).Value = 3489.1667
$ws.Cells.Item(105, 13).Value = -1742.1667

$ws.Cells.Item(134, 8).Value = 4939.75
$ws.Cells.Item(134, 9).Value = 5383.524
$ws.Cells.Item(134, 11).Value = 16150.572
$ws.Cells.Item(134, 13).Value = -13615.572

$ws.Cells.Item(140, 8).Value = 39500
$ws.Cells.Item(140, 10).Value = 39500
$ws.Cells.Item(140, 12).Value = 39500
$ws.Cells.Item(140, 14).Value = -49860

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1322
$ws.Cells.Item(16, 9).Value = 1340
$ws.Cells.Item(16, 10).Value = 1250
$ws.Cells.Item(16, 11).Value = 1340
$ws.Cells.Item(16, 12).Value = 1250
$ws.Cells.Item(16, 13).Value = -1053
$ws.Cells.Item(16, 14).Value = -1824

$ws.Cells.Item(22, 8).Value = 429.125
$ws.Cells.Item(22, 9).Value = 110.333336
$ws.Cells.Item(22, 11).Value = 110.333336
$ws.Cells.Item(22, 13).Value = 239.666664

$ws.Cells.Item(31, 8).Value = 14500.034
$ws.Cells.Item(31, 9).Value = 25714.846
$ws.Cells.Item(31, 10).Value = 5388
$ws.Cells.Item(31, 11).Value = 25714.846
$ws.Cells.Item(31, 12).Value = 5388
$ws.Cells.Item(31, 13).Value = -25419.846
$ws.Cells.Item(31, 14).Value = -5978

$ws.Cells.Item(34, 8).Value = 14500.034
$ws.Cells.Item(34, 9).Value = 25714.846
$ws.Cells.Item(34, 10).Value = 5388
$ws.Cells.Item(34, 11).Value = 25714.846
$ws.Cells.Item(34, 12).Value = 5388
$ws.Cells.Item(34, 13).Value = -25512.846
$ws.Cells.Item(34, 14).Value = -5792

$ws.Cells.Item(113, 8).Value = 1322
$ws.Cells.Item(113, 9).Value = 1340
$ws.Cells.Item(113, 10).Value = 1250
$ws.Cells.Item(113, 11).Value = 1340
$ws.Cells.Item(113, 12).Value = 1250
$ws.Cells.Item(113, 13).Value = 830
$ws.Cells.Item(113, 14).Value = -5590

$ws.Cells.Item(132, 8).Value = 11485.642
$ws.Cells.Item(132, 9).Value = 15532.723
$ws.Cells.Item(132, 10).Value = 2915.353
$ws.Cells.Item(132, 11).Value = 46598.169
$ws.Cells.Item(132, 12).Value = 8746.059000000001
$ws.Cells.Item(132, 13).Value = -44068.169
$ws.Cells.Item(132, 14).Value = -13806.059

$ws.Cells.Item(134, 8).Value = 1479.2222
$ws.Cells.Item(134, 9).Value = 1383.1666
$ws.Cells.Item(134, 10).Value = 1671.3334
$ws.Cells.Item(134, 11).Value = 4149.4998
$ws.Cells.Item(134, 12).Value = 5014.0002
$ws.Cells.Item(134, 13).Value = -1614.4998
$ws.Cells.Item(134, 14).Value = -10084.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 4437.9165
$ws.Cells.Item(107, 10).Value = 284.73685
$ws.Cells.Item(107, 12).Value = 854.21055
$ws.Cells.Item(107, 14).Value = -4694.21055

$ws.Cells.Item(129, 8).Value = 1299.9286
$ws.Cells.Item(129, 10).Value = 1387.9
$ws.Cells.Item(129, 12).Value = 4163.700000000001
$ws.Cells.Item(129, 14).Value = -14163.7

$ws.Cells.Item(131, 8).Value = 758.13
$ws.Cells.Item(131, 10).Value = 779.2083
$ws.Cells.Item(131, 12).Value = 2337.6249
$ws.Cells.Item(131, 14).Value = -12417.6249

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 16672500
$ws.Cells.Item(52, 10).Value = 16672500
$ws.Cells.Item(52, 12).Value = 16672500
$ws.Cells.Item(52, 14).Value = -16673018

$ws.Cells.Item(132, 8).Value = 17272.223
$ws.Cells.Item(132, 9).Value = 3743.1538
$ws.Cells.Item(132, 10).Value = 52447.8
$ws.Cells.Item(132, 11).Value = 11229.4614
$ws.Cells.Item(132, 12).Value = 157343.4
$ws.Cells.Item(132, 13).Value = -8699.4614
$ws.Cells.Item(132, 14).Value = -162403.4

$ws.Cells.Item(135, 8).Value = 50390
$ws.Cells.Item(135, 10).Value = 50390
$ws.Cells.Item(135, 12).Value = 50390
$ws.Cells.Item(135, 14).Value = -60530

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 2779.0588
$ws.Cells.Item(100, 9).Value = 1114.6666
$ws.Cells.Item(100, 10).Value = 3135.7144
$ws.Cells.Item(100, 11).Value = 1114.6666
$ws.Cells.Item(100, 12).Value = 3135.7144
$ws.Cells.Item(100, 13).Value = -573.6666
$ws.Cells.Item(100, 14).Value = -4217.7144

$ws.Cells.Item(104, 8).Value = 28087.8
$ws.Cells.Item(104, 10).Value = 28087.8
$ws.Cells.Item(104, 12).Value = 28087.8
$ws.Cells.Item(104, 14).Value = -35075.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 3299.8
$ws.Cells.Item(15, 10).Value = 3299.8
$ws.Cells.Item(15, 12).Value = 3299.8
$ws.Cells.Item(15, 14).Value = -3875.8

$ws.Cells.Item(81, 8).Value = 125001150
$ws.Cells.Item(81, 9).Value = 1446.4
$ws.Cells.Item(81, 11).Value = 2892.8
$ws.Cells.Item(81, 13).Value = -1831.8

$ws.Cells.Item(84, 8).Value = 125001150
$ws.Cells.Item(84, 9).Value = 1446.4
$ws.Cells.Item(84, 11).Value = 14464
$ws.Cells.Item(84, 13).Value = -9160

$ws.Cells.Item(113, 8).Value = 1288031
$ws.Cells.Item(113, 9).Value = 1614.75
$ws.Cells.Item(113, 10).Value = 3003252.5
$ws.Cells.Item(113, 11).Value = 4844.25
$ws.Cells.Item(113, 12).Value = 9009757.5
$ws.Cells.Item(113, 13).Value = -2674.25
$ws.Cells.Item(113, 14).Value = -9014097.5
